$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 165, shifting rows 165:236 down to 166:237.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new data record (copy the
# unchanged columns from the row that is now at 166, then overwrite the
# columns that actually differ).
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = 44636
$ws.Cells.Item(165, 4).NumberFormat = $ws.Cells.Item(166, 4).NumberFormat
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100102
$ws.Cells.Item(165, 8).Value = "Cítricos"
$ws.Cells.Item(165, 9).Value = 100102005
$ws.Cells.Item(165, 10).Value = "Naranja"
$ws.Cells.Item(165, 11).Value = "Valencia"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 220
$ws.Cells.Item(165, 14).Value = 9000
$ws.Cells.Item(165, 15).Value = 10000
$ws.Cells.Item(165, 16).Value = 9545
$ws.Cells.Item(165, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(165, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value = 636
$ws.Cells.Item(165, 20).Value = 15
